$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ID" row (old row 2) - rows shift up by one
$ws.Rows(2).Delete()

# Prefab / NormalStateFunc / UpStateFunc / Desc rows change their Type column from "int" to "string"
$ws.Range("B4").Value = "string"
$ws.Range("B5").Value = "string"
$ws.Range("B6").Value = "string"
$ws.Range("B7").Value = "string"

# Fix up the data validation range that the row deletion shrank incorrectly
$ws.Range("F8:F1048576").Validation.Delete()
$ws.Range("F8:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Update the remembered selection
$ws.Range("G14").Select()
